$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Naive QoQ error series: a new (most-recent) error value is inserted at the
# front of each data row (column B), pushing the existing B:J values one
# column to the right (C:K). Any value that would spill past column K (the
# oldest tracked error) is dropped, since the window only keeps 10 columns.
$newFirstValues = @{
    2  = -0.5825945370336409
    3  = 0.09567504080935779
    4  = -0.2604190369987228
    5  = 0.8354549961584912
    6  = -0.1000793599026215
    7  = -0.3537865060796963
    8  = 0.1481773904324453
    9  = 0.157445989004155
    10 = -0.5006594565260708
    11 = 0.2803578805354692
    12 = -0.1719748578450117
    13 = 0.3058625397463315
    14 = -0.6123299526872862
    15 = 0.6883713851991116
    16 = -0.2766911554241067
}

$lastCol = 11   # column K
$firstCol = 2   # column B

foreach ($row in $newFirstValues.Keys) {
    # Read the existing values in this row (B..K) before overwriting them.
    $existing = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($row, $c)
        if ($cell.Value2 -ne $null) {
            $existing += $cell.Value2
        }
    }

    # Build the shifted row: new value first, then the old values, capped at
    # the number of columns available (B..K = 10 columns).
    $shifted = @($newFirstValues[$row]) + $existing
    $maxCount = $lastCol - $firstCol + 1
    if ($shifted.Length -gt $maxCount) {
        $shifted = $shifted[0..($maxCount - 1)]
    }

    for ($i = 0; $i -lt $shifted.Length; $i++) {
        $ws.Cells.Item($row, $firstCol + $i).Value = $shifted[$i]
    }
}
